$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Locate the target paragraph: "Change color of PUA by button - Laura"
#     (the last paragraph in the "Change color of PUA by button" task item,
#     right before the final PUA-import related bullets).
$hostPara = $null
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "Change color of PUA by button - Laura*") {
        $hostPara = $candidate
    }
}

# --- Part 1: split the trailing " - Laura" run into three runs:
#     " " (keeps the original run's rsid), an en dash "\u2013", and " Laura".
$paraRange = $hostPara.Range.Duplicate
$searchRange = $hostPara.Range.Duplicate
$findOk = $searchRange.Find.Execute(" - Laura", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$leadingText = $d.Range($paraRange.Start, $searchRange.Start).Text

$newRunsXml = "<w:r><w:rPr><w:lang w:val=""en-US""/></w:rPr><w:t xml:space=""preserve"">" + $leadingText + "</w:t></w:r>" `
    + "<w:r w:rsidR=""00E4227F""><w:rPr><w:lang w:val=""en-US""/></w:rPr><w:t xml:space=""preserve""> </w:t></w:r>" `
    + "<w:r><w:rPr><w:lang w:val=""en-US""/></w:rPr><w:t>&#8211;</w:t></w:r>" `
    + "<w:r><w:rPr><w:lang w:val=""en-US""/></w:rPr><w:t xml:space=""preserve""> Laura</w:t></w:r>"
$paraContentXml = "<w:p $wNs>$newRunsXml</w:p>"

$contentRange = $d.Range($paraRange.Start, $paraRange.End)
$contentRange.InsertXML($paraContentXml)

# --- Part 2: insert a new, empty list paragraph (same list style/level)
#     immediately after the paragraph that was just edited.
$hostPara = $d.Paragraphs($hostPara.Index)
$insertionPoint = $d.Range($hostPara.Range.End, $hostPara.Range.End)
$newParaXml = "<w:p $wNs><w:pPr><w:pStyle w:val=""Listenabsatz""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr><w:rPr><w:lang w:val=""en-US""/></w:rPr></w:pPr></w:p>"
$insertionPoint.InsertXML($newParaXml)
